# Weekly refresh of the "Achicoria" price series: drop the oldest
# record (row 36) by shifting the whole block of rows up by one,
# then append the newest week's record in the now-vacant last row (71).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 37-71 up into rows 36-70 (this overwrites row 36's old
# data and duplicates old row 71's data into new row 71).
$src = $ws.Range("A37:R71")
$dst = $ws.Range("A36:R70")
$src.Copy($dst)

# Update the newly appended last row (71) with this week's figures:
# new date and new volume; price columns stay the same as before.
$ws.Range("D71").Value = 45209
$ws.Range("J71").Value = 70
